$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 3 "Function Scope" -> Strict mode column: keep text "undefined", color it red
$cell = $t.Cell(3, 3)
$cell.Range.Font.Color = 255

# Row 5 "Object not directly*" -> Non-strict mode column: text changes to "Global object", colored red
$cell = $t.Cell(5, 2)
$cell.Range.Text = "Global object"
$cell.Range.Font.Color = 255

# Row 5 "Object not directly*" -> Strict mode column: keep text "undefined", color it red
$cell = $t.Cell(5, 3)
$cell.Range.Font.Color = 255
